$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet for the "Valid login" test case data
$ws.Name = "ValidLogin"

# Shift existing login data right/down and add header row + new column
# Final layout:
#   A1 = UserName   B1 = Password
#   A2 = admin      B2 = manager
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Match the recorded view state: zoom level and active selection
$ws.Range("B3").Select()
$excel.ActiveWindow.Zoom = 175
